# "last revisions of the final project"
#
# ID_4532 sheet: rename the "Weights of the days" label (B5) to the new
# text "Function weight days [0-1]" (this mints a new shared-string entry),
# widen column B to fit the new label, update the raw EMD input values in
# row 7 (6/12/18-month table) and row 18 (12-month patient table), and
# leave the selection on E24 as the last user action.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ID_4532")

# --- relabel B5 ---------------------------------------------------------
$ws.Range("B5").Value = "Function weight days [0-1]"

# widen column B so the longer label keeps fitting (22.42578125 -> 24.5703125)
$ws.Columns("B").ColumnWidth = 24.5703125

# --- updated EMD raw inputs (first table, row 7) ------------------------
$ws.Range("D7").Value = 439.39234016021902
$ws.Range("E7").Value = 654.28872805784101
$ws.Range("F7").Value = 702.68683785809503

# --- updated EMD raw inputs (second table, row 18) -----------------------
$ws.Range("C18").Value = 70.436045094786607
$ws.Range("D18").Value = 74.826132426441106
$ws.Range("E18").Value = 114.872977309681
$ws.Range("F18").Value = 150.17002954437899

# last thing the author did before saving: leave the selection on E24
$ws.Range("E24").Select()
